$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C2:C5) from 2023-10-25 (45224) to 2023-11-03 (45233)
$ws.Range("C2:C5").Value = 45233
